# Apply crypto price/volume updates to match the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells whose new values would otherwise be
# auto-parsed by Excel as numbers (losing the original inline-string/text type).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the updated Coin / Link / Price / Volume(1h) values.
$ws.Range("D2").Value = "69.427.85"
$ws.Range("E2").Value = "  -3.86%  "
$ws.Range("D3").Value = "2.506.80"
$ws.Range("E3").Value = "  -4.77%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "576.35"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "166.21"
$ws.Range("E6").Value = "  -4.39%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "2.505.51"
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("E10").Value = "  -6.43%  "
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "2.939.05"
$ws.Range("E14").Value = "  -5.62%  "
$ws.Range("D15").Value = "69.376.53"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("D17").Value = "24.81"
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("D18").Value = "2.501.12"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").Value = "7.80"
$ws.Range("E19").Value = "  -5.17%  "
$ws.Range("D20").Value = "11.34"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").Value = "347.34"
$ws.Range("E21").Value = "  -6.97%  "
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  -5.82%  "
$ws.Range("D25").Value = "68.51"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E26").Value = "  -5.84%  "
$ws.Range("D27").Value = "8.88"
$ws.Range("E27").Value = "  -6.84%  "
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -4.87%  "
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "461.97"
$ws.Range("E33").Value = "  -6.60%  "
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("D37").Value = "154.18"
$ws.Range("E37").Value = "  -5.45%  "
$ws.Range("D38").Value = "18.95"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("E44").Value = "  -14.15%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "38.10"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.29"
$ws.Range("E46").Value = "  -10.04%  "
$ws.Range("D47").Value = "143.16"
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").Value = "3.49"
$ws.Range("E49").Value = "  -3.89%  "
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("D51").Value = "0.0730"
$ws.Range("E51").Value = "  -1.54%  "
